# edit.ps1 - reproduces the text content changes described by the diff.
$p = $ppt.ActivePresentation

# --- Slide 18: "Batch Stages - Example" content placeholder ---
# Merge/retype the Job Name / Batch Stage lines (runs collapse, tabs/spacing
# change) while leaving the two descriptive bullet paragraphs' text intact.
$s18 = $p.Slides.Item(18)
$sh18 = $s18.Shapes.Item(2)
$tr18 = $sh18.TextFrame.TextRange

# Edit from the end of the text backwards so earlier character offsets stay valid.
$tr18.Characters(192, 24).Text = "Batch Stage`t:D000"
$tr18.Characters(161, 30).Text = "Job Name`t      :B.DATE.CHANGE"
$tr18.Characters(32, 20).Text = "`tBatch Stage`t: A000"
$tr18.Characters(1, 30).Text = "Job Name`t`t:EB.CYCLE.DATES"

# --- Slide 47: "Common Variable file Inserted" textbox ---
# Split the single run into four runs (identical formatting) so each word
# is its own run, matching the retyped source.
$s47 = $p.Slides.Item(47)
$sh47 = $s47.Shapes.Item(5)
$tr47 = $sh47.TextFrame.TextRange

$tr47.Characters(1, 7).Text = "Common "
$tr47.Characters(8, 9).Text = "Variable "
$tr47.Characters(17, 5).Text = "file "
$tr47.Characters(22, 8).Text = "Inserted"
